$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 842; this shifts the existing rows 842:883 down to 843:884
# (matching the dimension growing from A1:D883 to A1:D884).
$ws.Rows("842:842").Insert()

# Populate the newly inserted row 842 with the new data point:
# 2026/02/23, 月, 19, 26
# Column A holds a date-like string ("2026/02/23") that must stay plain text
# (matching the rest of the sheet, which stores dates as literal strings, not
# real dates). Temporarily mark the cell as Text before assigning the value so
# it isn't auto-converted into a date serial, then drop the explicit format so
# the cell keeps the sheet's default (unstyled) look.
$ws.Range("A842").NumberFormat = "@"
$ws.Range("A842").Value = "2026/02/23"
$ws.Range("A842").ClearFormats()

$ws.Range("B842").Value = "月"
$ws.Range("C842").Value = 19
$ws.Range("D842").Value = 26
